$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 188, shifting existing rows 188-274 down to 189-275.
$ws.Rows.Item(188).Insert()

# Populate the newly inserted row 188 with the new weekly price record.
$ws.Cells.Item(188, 1).Value = 5
$ws.Cells.Item(188, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(188, 3).Value = "Maule"
$ws.Cells.Item(188, 4).Value = (Get-Date -Year 2022 -Month 2 -Day 17 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(188, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(188, 5).Value = 7
$ws.Cells.Item(188, 6).Value = 100112023
$ws.Cells.Item(188, 7).Value = "Brócoli"
$ws.Cells.Item(188, 8).Value = "Sin especificar"
$ws.Cells.Item(188, 9).Value = "Segunda"
$ws.Cells.Item(188, 10).Value = 3000
$ws.Cells.Item(188, 11).Value = 800
$ws.Cells.Item(188, 12).Value = 800
$ws.Cells.Item(188, 13).Value = 800
$ws.Cells.Item(188, 14).Value = "`$/unidad"
$ws.Cells.Item(188, 15).Value = "Región del Maule"
$ws.Cells.Item(188, 16).Value = 800
$ws.Cells.Item(188, 17).Value = 1
$ws.Cells.Item(188, 18).Value = "Hortaliza"
